# Apply "repull data, push all data, mean calculation" update to the dSF
# column (column F) on Sheet1. These values are re-pulled source data, so
# they are written directly as literal numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -1
    9  = -4
    10 = 4
    12 = 2
    13 = 3
    14 = 6
    17 = 1
    22 = 5
    23 = -1
    24 = -2
    30 = 2
    31 = 0
    35 = 1
    38 = 0
    60 = 2
    68 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
